$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70, shifting existing rows 70:141 down to 71:142
$ws.Rows.Item(70).Insert()

# Populate the new row 70 with the new data record
$ws.Range("A70").Value = 5
$ws.Range("B70").Value = "Macroferia Regional de Talca"
$ws.Range("C70").Value = "Maule"
$ws.Range("D70").Value = 44589
$ws.Range("E70").Value = 7
$ws.Range("F70").Value = 100112031
$ws.Range("G70").Value = "Poroto verde"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 150
$ws.Range("K70").Value = 30000
$ws.Range("L70").Value = 30000
$ws.Range("M70").Value = 30000
$ws.Range("N70").Value = "`$/saco 25 kilos"
$ws.Range("O70").Value = "Región del Maule"
$ws.Range("P70").Value = 1200
$ws.Range("Q70").Value = 25
$ws.Range("R70").Value = "Hortaliza"
